# Auto-generated edit script applying numeric corrections to H:N columns
# across multiple worksheets, per the commit diff (scheduled runner profit recompute).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value2 = 8374
$ws.Range("I86").Value2 = 8665.333000000001
$ws.Range("J86").Value2 = 7500
$ws.Range("K86").Value2 = 8665.333000000001
$ws.Range("L86").Value2 = 7500
$ws.Range("M86").Value2 = -7542.333000000001

$ws.Range("H89").Value2 = 8374
$ws.Range("I89").Value2 = 8665.333000000001
$ws.Range("J89").Value2 = 7500
$ws.Range("K89").Value2 = 43326.665
$ws.Range("L89").Value2 = 37500
$ws.Range("M89").Value2 = -37710.665

$ws.Range("H137").Value2 = 62505212
$ws.Range("I137").Value2 = 50003908
$ws.Range("J137").Value2 = 83340720
$ws.Range("K137").Value2 = 150011724
$ws.Range("L137").Value2 = 250022160
$ws.Range("M137").Value2 = -150009174
$ws.Range("N137").Value2 = -250027260

$ws.Range("H138").Value2 = 4338122.5
$ws.Range("I138").Value2 = 4935.75
$ws.Range("J138").Value2 = 4840521
$ws.Range("K138").Value2 = 14807.25
$ws.Range("L138").Value2 = 14521563
$ws.Range("M138").Value2 = -9667.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 18879412
$ws.Range("I32").Value2 = 21285254
$ws.Range("J32").Value2 = 33654.168
$ws.Range("K32").Value2 = 21285254
$ws.Range("L32").Value2 = 33654.168
$ws.Range("M32").Value2 = -21284967

$ws.Range("H45").Value2 = 2338
$ws.Range("I45").Value2 = 2199.8572
$ws.Range("J45").Value2 = 2579.75
$ws.Range("K45").Value2 = 2199.8572
$ws.Range("L45").Value2 = 2579.75
$ws.Range("M45").Value2 = -1822.8572

$ws.Range("H74").Value2 = 66668400
$ws.Range("I74").Value2 = 83334830
$ws.Range("J74").Value2 = 2672.8333
$ws.Range("K74").Value2 = 83334830
$ws.Range("L74").Value2 = 2672.8333
$ws.Range("M74").Value2 = -83333956

$ws.Range("H77").Value2 = 66668400
$ws.Range("I77").Value2 = 83334830
$ws.Range("J77").Value2 = 2672.8333
$ws.Range("K77").Value2 = 416674150
$ws.Range("L77").Value2 = 13364.1665
$ws.Range("M77").Value2 = -416669782

$ws.Range("H132").Value2 = 33343382
$ws.Range("I132").Value2 = 11576.28
$ws.Range("J132").Value2 = 200002400
$ws.Range("K132").Value2 = 34728.84
$ws.Range("L132").Value2 = 600007200
$ws.Range("M132").Value2 = -32198.84
$ws.Range("N132").Value2 = -600012260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value2 = 119999.5
$ws.Range("I59").Value2 = 0
$ws.Range("J59").Value2 = 119999.5
$ws.Range("K59").Value2 = 0
$ws.Range("L59").Value2 = 119999.5
$ws.Range("N59").Value2 = -121693.5

$ws.Range("H74").Value2 = 86997.5
$ws.Range("I74").Value2 = 73995
$ws.Range("J74").Value2 = 100000
$ws.Range("K74").Value2 = 73995
$ws.Range("L74").Value2 = 100000
$ws.Range("M74").Value2 = -73059

$ws.Range("H77").Value2 = 86997.5
$ws.Range("I77").Value2 = 73995
$ws.Range("J77").Value2 = 100000
$ws.Range("K77").Value2 = 221985
$ws.Range("L77").Value2 = 300000
$ws.Range("M77").Value2 = -217305

$ws.Range("H134").Value2 = 2276.0588
$ws.Range("I134").Value2 = 2399.7144
$ws.Range("J134").Value2 = 1699
$ws.Range("K134").Value2 = 7199.1432
$ws.Range("L134").Value2 = 5097
$ws.Range("M134").Value2 = -4664.1432
$ws.Range("N134").Value2 = -10167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value2 = 800
$ws.Range("I14").Value2 = 0
$ws.Range("J14").Value2 = 800
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = 800
$ws.Range("N14").Value2 = -1140

$ws.Range("H107").Value2 = 2765.111
$ws.Range("I107").Value2 = 2557.2856
$ws.Range("J107").Value2 = 3492.5
$ws.Range("K107").Value2 = 2557.2856
$ws.Range("L107").Value2 = 3492.5
$ws.Range("M107").Value2 = -637.2856000000002

$ws.Range("H132").Value2 = 77295.664
$ws.Range("I132").Value2 = 89345.30499999999
$ws.Range("J132").Value2 = 8010.25
$ws.Range("K132").Value2 = 268035.915
$ws.Range("L132").Value2 = 24030.75
$ws.Range("M132").Value2 = -265505.915

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 4742.5938
$ws.Range("I5").Value2 = 669
$ws.Range("J5").Value2 = 11531.917
$ws.Range("K5").Value2 = 2007
$ws.Range("L5").Value2 = 34595.751
$ws.Range("M5").Value2 = -1895

$ws.Range("H32").Value2 = 92658.91
$ws.Range("I32").Value2 = 333666.66
$ws.Range("J32").Value2 = 2281
$ws.Range("K32").Value2 = 1000999.98
$ws.Range("L32").Value2 = 6843
$ws.Range("M32").Value2 = -1000716.98
$ws.Range("N32").Value2 = -7409

$ws.Range("H113").Value2 = 0
$ws.Range("I113").Value2 = 0
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 0
$ws.Range("L113").Value2 = 0
$ws.Range("N113").ClearContents()

$ws.Range("H131").Value2 = 1893.0769
$ws.Range("I131").Value2 = 1926.6666
$ws.Range("J131").Value2 = 1883
$ws.Range("K131").Value2 = 5779.9998
$ws.Range("L131").Value2 = 5649
$ws.Range("M131").Value2 = -739.9997999999996
$ws.Range("N131").Value2 = -15729

$ws.Range("H135").Value2 = 4742.5938
$ws.Range("I135").Value2 = 669
$ws.Range("J135").Value2 = 11531.917
$ws.Range("K135").Value2 = 6021
$ws.Range("L135").Value2 = 103787.253
$ws.Range("M135").Value2 = -3486

$ws.Range("H137").Value2 = 3199.6667
$ws.Range("I137").Value2 = 3199.6667
$ws.Range("J137").Value2 = 0
$ws.Range("K137").Value2 = 9599.000100000001
$ws.Range("L137").Value2 = 0
$ws.Range("M137").Value2 = -4499.000100000001
$ws.Range("N137").ClearContents()

$ws.Range("H140").Value2 = 2638.1785
$ws.Range("I140").Value2 = 2660.6667
$ws.Range("J140").Value2 = 2612.2307
$ws.Range("K140").Value2 = 7982.000100000001
$ws.Range("L140").Value2 = 7836.6921
$ws.Range("M140").Value2 = -2802.000100000001
$ws.Range("N140").Value2 = -18196.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 5492.0454
$ws.Range("I132").Value2 = 4801.9033
$ws.Range("J132").Value2 = 7137.769
$ws.Range("K132").Value2 = 14405.7099
$ws.Range("L132").Value2 = 21413.307
$ws.Range("M132").Value2 = -11875.7099
$ws.Range("N132").Value2 = -26473.307

$ws.Range("H133").Value2 = 149999.5
$ws.Range("I133").Value2 = 0
$ws.Range("J133").Value2 = 149999.5
$ws.Range("K133").Value2 = 0
$ws.Range("L133").Value2 = 149999.5
$ws.Range("N133").Value2 = -160119.5

$ws.Range("H137").Value2 = 150000
$ws.Range("I137").Value2 = 0
$ws.Range("J137").Value2 = 150000
$ws.Range("K137").Value2 = 0
$ws.Range("L137").Value2 = 150000
$ws.Range("N137").Value2 = -160200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value2 = 0
$ws.Range("I19").Value2 = 0
$ws.Range("J19").Value2 = 0
$ws.Range("K19").Value2 = 0
$ws.Range("L19").Value2 = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value2 = 11498.5
$ws.Range("I18").Value2 = 0
$ws.Range("J18").Value2 = 11498.5
$ws.Range("K18").Value2 = 0
$ws.Range("L18").Value2 = 11498.5
$ws.Range("N18").Value2 = -11844.5

$ws.Range("H62").Value2 = 11225
$ws.Range("I62").Value2 = 12000
$ws.Range("J62").Value2 = 10966.667
$ws.Range("K62").Value2 = 12000
$ws.Range("L62").Value2 = 10966.667
$ws.Range("M62").Value2 = -11376
$ws.Range("N62").Value2 = -12214.667

$ws.Range("H65").Value2 = 11225
$ws.Range("I65").Value2 = 12000
$ws.Range("J65").Value2 = 10966.667
$ws.Range("K65").Value2 = 60000
$ws.Range("L65").Value2 = 54833.335
$ws.Range("M65").Value2 = -56880
$ws.Range("N65").Value2 = -61073.335

$ws.Range("H81").Value2 = 1334.5714
$ws.Range("I81").Value2 = 961.8946999999999
$ws.Range("J81").Value2 = 4875
$ws.Range("K81").Value2 = 1923.7894
$ws.Range("L81").Value2 = 9750
$ws.Range("M81").Value2 = -862.7893999999999
$ws.Range("N81").Value2 = -11872

$ws.Range("H84").Value2 = 1334.5714
$ws.Range("I84").Value2 = 961.8946999999999
$ws.Range("J84").Value2 = 4875
$ws.Range("K84").Value2 = 9618.947
$ws.Range("L84").Value2 = 48750
$ws.Range("M84").Value2 = -4314.947
$ws.Range("N84").Value2 = -59358

$ws.Range("H107").Value2 = 1000
$ws.Range("I107").Value2 = 0
$ws.Range("J107").Value2 = 1000
$ws.Range("K107").Value2 = 0
$ws.Range("L107").Value2 = 3000
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value2 = -6840

